$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'72.056.29"
$ws.Range("E2").Value = "'  +0.45%  "
$ws.Range("D3").Value = "'4.043.54"
$ws.Range("E3").Value = "'  +0.00%  "
$ws.Range("E4").Value = "'  -0.06%  "
$ws.Range("D5").Value = "'537.64"
$ws.Range("E5").Value = "'  +0.70%  "
$ws.Range("D6").Value = "'149.44"
$ws.Range("E6").Value = "'  -1.98%  "
$ws.Range("D7").Value = "'4.039.26"
$ws.Range("E7").Value = "'  +0.13%  "
$ws.Range("E8").Value = "'  +0.63%  "
$ws.Range("E9").Value = "'  -0.09%  "
$ws.Range("D10").Value = "'0.752"
$ws.Range("E10").Value = "'  -1.07%  "
$ws.Range("D11").Value = "'0.172"
$ws.Range("E11").Value = "'  -1.75%  "
$ws.Range("D12").Value = "'53.36"
$ws.Range("E12").Value = "'  +9.41%  "
$ws.Range("D13").Value = "'0.0000331"
$ws.Range("E13").Value = "'  -0.44%  "
$ws.Range("E14").Value = "'  -0.13%  "
$ws.Range("D15").Value = "'4.690.25"
$ws.Range("E15").Value = "'  -0.03%  "
$ws.Range("D16").Value = "'4.050.04"
$ws.Range("E16").Value = "'  +0.02%  "
$ws.Range("D17").Value = "'14.32"
$ws.Range("E17").Value = "'  -0.49%  "
$ws.Range("D18").Value = "'20.74"
$ws.Range("E18").Value = "'  -1.25%  "
$ws.Range("E19").Value = "'  -1.12%  "
$ws.Range("E20").Value = "'  -1.14%  "
$ws.Range("D21").Value = "'72.062.20"
$ws.Range("E21").Value = "'  +0.44%  "
$ws.Range("D22").Value = "'437.91"
$ws.Range("E22").Value = "'  +0.23%  "
$ws.Range("D23").Value = "'98.26"
$ws.Range("E23").Value = "'  -1.20%  "
$ws.Range("E24").Value = "'  -6.10%  "
$ws.Range("D25").Value = "'4.26"
$ws.Range("E25").Value = "'  +0.87%  "
$ws.Range("D26").Value = "'14.53"
$ws.Range("E26").Value = "'  -1.91%  "
$ws.Range("D27").Value = "'4.35"
$ws.Range("E27").Value = "'  +23.48%  "
$ws.Range("D28").Value = "'11.31"
$ws.Range("E28").Value = "'  -0.89%  "
$ws.Range("E29").Value = "'  -1.92%  "
$ws.Range("E30").Value = "'  +2.02%  "
$ws.Range("D31").Value = "'37.12"
$ws.Range("E31").Value = "'  -0.28%  "
$ws.Range("D32").Value = "'8.47"
$ws.Range("E32").Value = "'  +24.30%  "
$ws.Range("E33").Value = "'  +2.00%  "
$ws.Range("D34").Value = "'50.00"
$ws.Range("E34").Value = "'  +16.17%  "
$ws.Range("D35").Value = "'13.58"
$ws.Range("E35").Value = "'  -0.94%  "
$ws.Range("D36").Value = "'676.90"
$ws.Range("E36").Value = "'  -0.44%  "
$ws.Range("D37").Value = "'66.40"
$ws.Range("E37").Value = "'  -0.41%  "
$ws.Range("D38").Value = "'0.459"
$ws.Range("E38").Value = "'  +5.91%  "
$ws.Range("D39").Value = "'0.0₃0880"
$ws.Range("E39").Value = "'  +2.34%  "
$ws.Range("E40").Value = "'  +8.05%  "
$ws.Range("D41").Value = "'0.148"
$ws.Range("E41").Value = "'  -6.26%  "
$ws.Range("D42").Value = "'3.40"
$ws.Range("E42").Value = "'  -2.13%  "
$ws.Range("E43").Value = "'  +17.70%  "
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "'  +0.09%  "
$ws.Range("E45").Value = "'  -1.03%  "
$ws.Range("D46").Value = "'1.00"
$ws.Range("E46").Value = "'  +0.02%  "
$ws.Range("E47").Value = "'  -1.01%  "
$ws.Range("D48").Value = "'2.68"
$ws.Range("E48").Value = "'  -1.87%  "
$ws.Range("E49").Value = "'  +2.04%  "
$ws.Range("D50").Value = "'3.31"
$ws.Range("E50").Value = "'  -3.23%  "
$ws.Range("D51").Value = "'2.831.93"
$ws.Range("E51").Value = "'  +8.10%  "
